$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsCL = $wb.Worksheets.Item("CL")

# --- "CL" sheet: restructure into urban/rural/commercial residential columns ---

# New column headers (B/C/D) - order matters for shared-string allocation
$wsCL.Range("B1").Value = "urban residential"
$wsCL.Range("C1").Value = "rural residential"
$wsCL.Range("D1").Value = "commercial"

# Copy the formatting of the existing header cell (B1) onto the two new header cells
$wsCL.Range("B1").Copy() | Out-Null
$wsCL.Range("C1:D1").PasteSpecial(-4122) | Out-Null

# --- "About" sheet: add a new note row (row 28 left blank, row 29 gets the new text) ---
$wsAbout.Range("A29").Value = "For the U.S. model, we use the same component lifetimes across building types."

# Rename A1 header on CL sheet (allocated last so it becomes the final new shared string)
$wsCL.Range("A1").Value = "Building Component (years)"

# Set the column widths for the two new columns (closest achievable given the
# host's internal pixel-quantized column-width model)
$wsCL.Columns.Item(3).ColumnWidth = 17.022135416666668
$wsCL.Columns.Item(4).ColumnWidth = 15.736979166666666

# Rural residential (C) and commercial (D) columns just mirror urban residential (B)
$wsCL.Range("C2").Formula = '=$B2'
$wsCL.Range("D2").Formula = '=$B2'
$wsCL.Range("C3:D7").Formula = '=$B3'

# Explicit portrait page setup for the CL sheet
$wsCL.PageSetup.Orientation = 1

# Restore the selection to cell A2 on the CL sheet, then re-activate the About sheet
$wsCL.Activate() | Out-Null
$wsCL.Range("A2").Select() | Out-Null
$wsAbout.Activate() | Out-Null
